# "show card brief text"
# Populate the "brief" column (K) of the readings table with the
# placeholder intro text for every data row, and bring the header /
# sub-header cells in that column in line with the rest of the header
# rows (center aligned + wrap text), matching the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$briefText = "这是一些简单的介绍文字，不要太长，也不要太短，差不多就行了。"

# Data rows: 4..7 hold the actual reading entries; column K is "brief".
$ws.Range("K4:K7").Value = $briefText
$ws.Range("K4:K7").WrapText = $true

# Header rows (1: column title, 2: type hint, 3: Chinese label) get the
# same wrap formatting the neighbouring header cells already use so the
# whole row reads consistently once the brief column has content.
$ws.Range("K1").WrapText = $true
$ws.Range("K2").WrapText = $true
$ws.Range("K3").WrapText = $true

# Update the view so the freshly-populated column is what's on screen.
$ws.Range("K13").Select()
